$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.352.32"
$ws.Range("E2").Value = "  +1.76%  "

$ws.Range("D3").Value = "3.595.36"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'243.46"
$ws.Range("E5").Value = "  +2.41%  "

$ws.Range("D6").Value = "'1.81"
$ws.Range("E6").Value = "  +17.92%  "

$ws.Range("D7").Value = "'654.42"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +5.66%  "

$ws.Range("D9").Value = "'1.08"
$ws.Range("E9").Value = "  +4.80%  "

$ws.Range("D10").Value = "'1.00"
$ws.Range("E10").Value = "  -0.02%  "

$ws.Range("D11").Value = "3.592.52"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").Value = "'44.94"
$ws.Range("E12").Value = "  +4.82%  "

$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").Value = "'6.46"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").Value = "4.263.71"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").Value = "97.223.62"
$ws.Range("E16").Value = "  +1.80%  "

$ws.Range("E17").Value = "  +2.02%  "

$ws.Range("D18").Value = "3.581.42"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").Value = "'7.77"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").Value = "'12.62"
$ws.Range("E20").Value = "  -1.33%  "

$ws.Range("D21").Value = "'18.29"
$ws.Range("E21").Value = "  +1.83%  "

$ws.Range("D22").Value = "'0.548"
$ws.Range("E22").Value = "  +10.93%  "

$ws.Range("D23").Value = "'3.50"
$ws.Range("E23").Value = "  +1.42%  "

$ws.Range("D24").Value = "'517.89"
$ws.Range("E24").Value = "  +1.14%  "

$ws.Range("D25").Value = "'0.0000205"
$ws.Range("E25").Value = "  +3.27%  "

$ws.Range("D26").Value = "'6.99"
$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("D27").Value = "'103.22"
$ws.Range("E27").Value = "  +7.72%  "

$ws.Range("D28").Value = "'13.32"
$ws.Range("E28").Value = "  +4.01%  "

$ws.Range("E29").Value = "  +25.47%  "

$ws.Range("D30").Value = "'3.00"

$ws.Range("D31").Value = "'12.08"
$ws.Range("E31").Value = "  +4.47%  "

$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.18%  "

$ws.Range("D33").Value = "'0.190"
$ws.Range("E33").Value = "  +6.32%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "'31.92"
$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("E36").Value = "  +2.33%  "

$ws.Range("E37").Value = "  +3.72%  "

$ws.Range("D38").Value = "'618.54"
$ws.Range("E38").Value = "  +3.95%  "

$ws.Range("D39").Value = "'8.77"
$ws.Range("E39").Value = "  +2.36%  "

$ws.Range("E40").Value = "  +1.65%  "

$ws.Range("E41").Value = "  +2.81%  "

$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("E43").Value = "  +1.61%  "

$ws.Range("D44").Value = "'0.449"
$ws.Range("E44").Value = "  +39.08%  "

$ws.Range("D45").Value = "'6.14"
$ws.Range("E45").Value = "  +5.47%  "

$ws.Range("D46").Value = "'0.0454"
$ws.Range("E46").Value = "  +8.67%  "

$ws.Range("D47").Value = "'2.34"
$ws.Range("E47").Value = "  +1.09%  "

$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'8.79"
$ws.Range("E48").Value = "  +6.61%  "

$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "'23.66"
$ws.Range("E49").Value = "  +0.93%  "

$ws.Range("D50").Value = "'3.28"
$ws.Range("E50").Value = "  +6.59%  "

$ws.Range("D51").Value = "'32.71"
$ws.Range("E51").Value = "  -5.24%  "
